# Removing less than USD 5 price from extrapolation calibration because it is just a noise
# Updates recalculated values in columns D:H for rows 4, 9, 10, 11, 16, 18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 119110.6959536328
$ws.Range("E4").Value = -0.02540452173790457
$ws.Range("F4").Value = 0.2019830164175798
$ws.Range("G4").Value = -1.545695123913811
$ws.Range("H4").Value = 18.44252563185625

# Row 9
$ws.Range("D9").Value = 124921.4050677596
$ws.Range("E9").Value = -0.08349034835913985
$ws.Range("F9").Value = 0.3501939789064979
$ws.Range("G9").Value = -1.716480120115715
$ws.Range("H9").Value = 10.59627937947441

# Row 10
$ws.Range("D10").Value = 125804.6627585162
$ws.Range("E10").Value = -0.1238914850691796
$ws.Range("F10").Value = 0.4496105518733471
$ws.Range("G10").Value = -1.947917665978205
$ws.Range("H10").Value = 9.841464898599122

# Row 11
$ws.Range("D11").Value = 127900.3809220471
$ws.Range("E11").Value = -0.1361039053318512
$ws.Range("F11").Value = 0.4451099947975722
$ws.Range("G11").Value = -1.65527693074002
$ws.Range("H11").Value = 7.743714170899811

# Row 16
$ws.Range("D16").Value = 117583.6572386419
$ws.Range("E16").Value = -0.05578875267763517
$ws.Range("F16").Value = 0.1484356984337717
$ws.Range("G16").Value = -0.7013653358047669
$ws.Range("H16").Value = 9.826820736088662

# Row 18
$ws.Range("D18").Value = 118459.1535628688
$ws.Range("E18").Value = -0.01244252634384652
$ws.Range("F18").Value = 0.1565769665494396
$ws.Range("G18").Value = -0.18295335789966
$ws.Range("H18").Value = 5.932574827135706
